# Finalized changes to Front End SQL calls
#
# 1. Rename "Sheet1" -> "MENU BUTTONS"
# 2. Tidy up two existing SQL-call descriptions (replace enumerated id
#    lists with a compact "start-end" range notation).
# 3. Add the missing "Sql call" (column H) entries for the Cookie
#    Sandwich, French Fries, Fountain Drink and Ice cream cup rows,
#    matching the formatting already used by the neighbouring
#    "Description" cells in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the sheet ---------------------------------------------------
$ws.Name = "MENU BUTTONS"

# --- 2. Tweak existing SQL-call text ----------------------------------------
# Chicken Tenders row (row 9, column H)
$ws.Range("H9").Value = "(x,y can be 1026-1031 ) UPDATE entrees SET inventory = inventory - 1 where entrees.id in (1018,1020,1021,1022,1023,1024); Update sauces set inventory = inventory - 1 where sauces.id in(x,y); Update entrees set inventory = inventory - 3 where entrees.id = 1009;"

# Aggie Shakes row (row 11, column H)
$ws.Range("H11").Value = "(x can be 1032-1034) UPDATE desserts SET inventory = inventory - 1 where desserts.id in (x, 1036); UPDATE entrees SET inventory = inventory - 1 where entrees.id in (1023);"

# --- 3. Add the new "Sql call" entries --------------------------------------
# Cookie Sandwich row (row 12, column H) - previously empty
$ws.Range("H12").Value = "UPDATE desserts SET inventory = inventory - 2 where desserts.id in (1038); UPDATE desserts SET inventory = inventory - 1 where desserts.id in (1033, 1037); UPDATE entrees SET inventory = inventory - 1 where entrees.id in (1024);"

# French Fries row (row 13, column H) - previously empty
$ws.Range("H13").Value = "UPDATE entrees SET inventory = inventory - 2 where entrees.id in (1018, 1020, 1021, 1024);"

# Fountain Drink row (row 14, column H) - previously empty
$ws.Range("H14").Value = "(x can be 1000-1005) UPDATE drinks SET inventory = inventory - 1 where drinks.id in (x, 1006);"

# Ice cream cup row (row 15, column H) - previously empty
$ws.Range("H15").Value = "UPDATE desserts SET inventory = inventory - 2 where desserts.id in (x, 1037); UPDATE entrees SET inventory = inventory - 1 where entrees.id in (1023, 1024);"

# Match the new cells' formatting to the rest of row (same style as the
# "Description" cell in column G of each row).
foreach ($row in 12..15) {
    $ws.Range("G$row").Copy()
    $ws.Range("H$row").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
